$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '37.821.88'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '2.103.73'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '235.15'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').Value = '  +1.34%  '
$ws.Range('D7').Value = '58.24'
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.391'
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('D10').Value = '0.0779'
$ws.Range('E10').Value = '  +2.73%  '
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').Value = '2.415.52'
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('D13').Value = '14.48'
$ws.Range('E13').Value = '  +1.14%  '
$ws.Range('D14').Value = '21.29'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = '0.785'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '5.22'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('D17').Value = '2.108.73'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('D18').Value = '37.790.87'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').Value = '6.24'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').Value = '70.30'
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').Value = '0.0₃0823'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '227.47'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('D26').Value = '167.86'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('D27').Value = '8.96'
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  +3.65%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '1.41'
$ws.Range('E29').Value = '  -3.41%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '19.56'
$ws.Range('E30').Value = '  +2.63%  '
$ws.Range('E31').Value = '  +1.24%  '
$ws.Range('D32').Value = '4.62'
$ws.Range('E32').Value = '  +3.20%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = '2.60'
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.0623'
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').Value = '3.51'
$ws.Range('E36').Value = '  +7.05%  '
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = '5.43'
$ws.Range('E39').Value = '  -7.21%  '
$ws.Range('E40').Value = '  +3.43%  '
$ws.Range('D41').Value = '2.95'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.472.65'
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '96.69'
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('D45').Value = '1.17'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '4.15'
$ws.Range('E46').Value = '  -10.07%  '
$ws.Range('E47').Value = '  +2.95%  '
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('D50').Value = '3.02'
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('D51').Value = '2.299.50'
$ws.Range('E51').Value = '  +2.23%  '
